$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the glass-preference text "ALL_PREF_CDGM+SCHOTT_ONLY" with "SCHOTT"
# in the vendor-preference column (F) for the affected surface rows.
$ws.Range("F6").Value = "SCHOTT"
$ws.Range("F13").Value = "SCHOTT"
$ws.Range("F16").Value = "SCHOTT"
$ws.Range("F18").Value = "SCHOTT"
$ws.Range("F20").Value = "SCHOTT"

# Update the active selection on the sheet to N23.
$ws.Range("N23").Select()
